$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data rows 9..22 down to 10..23 (process bottom-up so we don't
# clobber a row before it has been copied to the next one).
for ($r = 22; $r -ge 9; $r--) {
    $target = $r + 1
    $ws.Cells.Item($target, 4).Value2  = $ws.Cells.Item($r, 4).Value2   # D Fecha
    $ws.Cells.Item($target, 10).Value2 = $ws.Cells.Item($r, 10).Value2  # J Volumen
    $ws.Cells.Item($target, 11).Value2 = $ws.Cells.Item($r, 11).Value2  # K Precio minimo
    $ws.Cells.Item($target, 12).Value2 = $ws.Cells.Item($r, 12).Value2  # L Precio maximo
    $ws.Cells.Item($target, 13).Value2 = $ws.Cells.Item($r, 13).Value2  # M Precio promedio ponderado
    $ws.Cells.Item($target, 15).Value2 = $ws.Cells.Item($r, 15).Value2  # O Origen
    $ws.Cells.Item($target, 16).Value2 = $ws.Cells.Item($r, 16).Value2  # P Precio $/Kg

    # Remaining columns are identical across all data rows, but copy them
    # too so row 23 is fully populated like the others.
    $ws.Cells.Item($target, 1).Value2  = $ws.Cells.Item($r, 1).Value2   # A Mercado ID
    $ws.Cells.Item($target, 2).Value2  = $ws.Cells.Item($r, 2).Value2   # B Mercado
    $ws.Cells.Item($target, 3).Value2  = $ws.Cells.Item($r, 3).Value2   # C Region
    $ws.Cells.Item($target, 5).Value2  = $ws.Cells.Item($r, 5).Value2   # E Codreg
    $ws.Cells.Item($target, 6).Value2  = $ws.Cells.Item($r, 6).Value2   # F Categoria ID
    $ws.Cells.Item($target, 7).Value2  = $ws.Cells.Item($r, 7).Value2   # G Categoria
    $ws.Cells.Item($target, 8).Value2  = $ws.Cells.Item($r, 8).Value2   # H Variedad
    $ws.Cells.Item($target, 9).Value2  = $ws.Cells.Item($r, 9).Value2   # I Calidad
    $ws.Cells.Item($target, 14).Value2 = $ws.Cells.Item($r, 14).Value2  # N Unidad de comercializacion
    $ws.Cells.Item($target, 17).Value2 = $ws.Cells.Item($r, 17).Value2  # Q Kg o Unidades
    $ws.Cells.Item($target, 18).Value2 = $ws.Cells.Item($r, 18).Value2  # R Clasificacion
}

# Row 23 is brand new; give its date cell the same number format used by
# the rest of column D (style index 2 / YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat

# New data for row 9
$ws.Cells.Item(9, 4).Value2  = 44799
$ws.Cells.Item(9, 10).Value2 = 500
$ws.Cells.Item(9, 11).Value2 = 10000
$ws.Cells.Item(9, 12).Value2 = 11000
$ws.Cells.Item(9, 13).Value2 = 10500
$ws.Cells.Item(9, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value2 = 420
